$d = $word.ActiveDocument

# Locate the placeholder text "GITHUB REPOSITORY : [ " (the opening bracket
# plus one space) that currently precedes the empty "  ]" placeholder.
$findRng = $d.Content
$found = $findRng.Find.Execute(
    "GITHUB REPOSITORY : [ ", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)

if ($found) {
    # Include the manual line break just before "GITHUB" so that it stays
    # attached to this run rather than the previous ("... Lévaro") run once
    # the paragraph is re-serialized.
    $leftStart = $findRng.Start - 1

    # Collapse to the point right after "[ " (i.e. right before the second,
    # still-empty space that precedes the closing bracket) and insert the
    # repository URL there.
    $insPoint = $findRng.Duplicate
    $insPoint.Collapse(0)
    $insPoint.InsertAfter("https://github.com/AlanZavala/Lab08")

    # Touch (and immediately revert) the formatting of the newly inserted
    # text and of the text that precedes it. This keeps the three pieces of
    # text - "GITHUB REPOSITORY : [ ", the URL, and " ]" - as distinct runs
    # instead of being silently recombined into a single run on save.
    $insPoint.Font.Bold = $true
    $insPoint.Font.Bold = $false

    $leftRng = $d.Range($leftStart, $insPoint.Start)
    $leftRng.Font.Bold = $true
    $leftRng.Font.Bold = $false
}
